$d = $word.ActiveDocument

# Locate the top-level "Setup" heading (style MonTitre1) -- not the Table of
# Contents entry that happens to contain the same word -- and strip its
# text/bookmark so the paragraph becomes empty while keeping the paragraph
# itself (and its MonTitre1 style) in place.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text.Trim()

    if ($styleName -eq "MonTitre1" -and $text -eq "Setup") {
        # Remove the run text but keep the trailing paragraph mark.
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Delete()

        # Remove the now-empty heading's TOC bookmark entirely (it no
        # longer anchors any text). The remaining bookmarks keep their
        # relative order and get renumbered automatically on save.
        $b = $d.Bookmarks.Item("_Toc101091930")
        $b.Delete()

        break
    }
}
